$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 36 (FY_4.png) - rows 36-44 (TO_1..TD_4) shift down to 37-45
$ws.Rows.Item(36).Insert()

# Apply updated values (recomputed metrics) across the dataset
$ws.Cells.Item(1, 2).Value = 0.473
$ws.Cells.Item(2, 2).Value = 0.974
$ws.Cells.Item(2, 3).Value = 0.003
$ws.Cells.Item(2, 4).Value = 0.398
$ws.Cells.Item(3, 2).Value = 0.79
$ws.Cells.Item(3, 3).Value = 0.003
$ws.Cells.Item(4, 2).Value = 0.441
$ws.Cells.Item(5, 2).Value = 0.336
$ws.Cells.Item(6, 2).Value = 1.271
$ws.Cells.Item(6, 3).Value = 0.004
$ws.Cells.Item(6, 4).Value = 0.118
$ws.Cells.Item(6, 5).Value = 0.4
$ws.Cells.Item(7, 2).Value = 0.977
$ws.Cells.Item(8, 2).Value = 1.269
$ws.Cells.Item(8, 3).Value = 0.004
$ws.Cells.Item(8, 4).Value = 0.645
$ws.Cells.Item(9, 2).Value = 1.164
$ws.Cells.Item(9, 3).Value = 0.004
$ws.Cells.Item(9, 4).Value = 0.225
$ws.Cells.Item(9, 5).Value = 0.9330000000000001
$ws.Cells.Item(10, 2).Value = 1.112
$ws.Cells.Item(10, 3).Value = 0.004
$ws.Cells.Item(10, 4).Value = 0.279
$ws.Cells.Item(11, 2).Value = 0.768
$ws.Cells.Item(11, 3).Value = 0.003
$ws.Cells.Item(12, 2).Value = 1.615
$ws.Cells.Item(12, 3).Value = 0.005
$ws.Cells.Item(12, 4).Value = 0.6860000000000001
$ws.Cells.Item(13, 2).Value = 0.734
$ws.Cells.Item(14, 2).Value = 0.792
$ws.Cells.Item(14, 3).Value = 0.003
$ws.Cells.Item(15, 2).Value = 0.831
$ws.Cells.Item(15, 3).Value = 0.003
$ws.Cells.Item(15, 4).Value = 0.412
$ws.Cells.Item(15, 5).Value = 0.667
$ws.Cells.Item(15, 6).Value = "Toni Ismail"
$ws.Cells.Item(15, 7).Value = "Benar"
$ws.Cells.Item(16, 2).Value = 0.6879999999999999
$ws.Cells.Item(16, 4).Value = 0.676
$ws.Cells.Item(17, 2).Value = 0.5679999999999999
$ws.Cells.Item(17, 3).Value = 0.002
$ws.Cells.Item(17, 4).Value = 0.383
$ws.Cells.Item(18, 2).Value = 0.854
$ws.Cells.Item(18, 4).Value = 0.45
$ws.Cells.Item(18, 5).Value = 0.667
$ws.Cells.Item(19, 2).Value = 0.511
$ws.Cells.Item(19, 3).Value = 0.002
$ws.Cells.Item(19, 4).Value = 0.445
$ws.Cells.Item(20, 2).Value = 0.914
$ws.Cells.Item(20, 3).Value = 0.003
$ws.Cells.Item(20, 4).Value = 0.297
$ws.Cells.Item(20, 5).Value = 0.867
$ws.Cells.Item(21, 2).Value = 0.385
$ws.Cells.Item(21, 4).Value = 0.343
$ws.Cells.Item(22, 2).Value = 0.946
$ws.Cells.Item(22, 4).Value = 0.171
$ws.Cells.Item(23, 2).Value = 0.901
$ws.Cells.Item(23, 3).Value = 0.003
$ws.Cells.Item(23, 4).Value = 0.378
$ws.Cells.Item(23, 5).Value = 1
$ws.Cells.Item(24, 2).Value = 1.793
$ws.Cells.Item(24, 3).Value = 0.005
$ws.Cells.Item(24, 4).Value = 0.515
$ws.Cells.Item(25, 2).Value = 1.266
$ws.Cells.Item(25, 3).Value = 0.004
$ws.Cells.Item(25, 4).Value = 0.545
$ws.Cells.Item(26, 2).Value = 1.143
$ws.Cells.Item(26, 3).Value = 0.004
$ws.Cells.Item(26, 4).Value = 0.094
$ws.Cells.Item(26, 5).Value = 0.867
$ws.Cells.Item(27, 2).Value = 1.028
$ws.Cells.Item(27, 4).Value = 0.546
$ws.Cells.Item(28, 2).Value = 1.045
$ws.Cells.Item(28, 4).Value = 0.569
$ws.Cells.Item(29, 2).Value = 0.662
$ws.Cells.Item(30, 2).Value = 1.139
$ws.Cells.Item(30, 3).Value = 0.004
$ws.Cells.Item(30, 4).Value = 0.534
$ws.Cells.Item(31, 2).Value = 0.7
$ws.Cells.Item(31, 3).Value = 0.002
$ws.Cells.Item(31, 4).Value = 0.275
$ws.Cells.Item(31, 5).Value = 0.9330000000000001
$ws.Cells.Item(32, 2).Value = 0.576
$ws.Cells.Item(32, 3).Value = 0.002
$ws.Cells.Item(32, 4).Value = 0.188
$ws.Cells.Item(33, 2).Value = 1.012
$ws.Cells.Item(33, 3).Value = 0.003
$ws.Cells.Item(33, 4).Value = 0.319
$ws.Cells.Item(33, 5).Value = 0.733
$ws.Cells.Item(33, 6).Value = "Fanny Yusuf"
$ws.Cells.Item(33, 7).Value = "Benar"
$ws.Cells.Item(34, 2).Value = 1.241
$ws.Cells.Item(34, 3).Value = 0.004
$ws.Cells.Item(34, 4).Value = 0.497
$ws.Cells.Item(34, 5).Value = 0.6
$ws.Cells.Item(35, 2).Value = 1.324
$ws.Cells.Item(35, 3).Value = 0.004
$ws.Cells.Item(35, 4).Value = 0.484
$ws.Cells.Item(35, 5).Value = 0.667
$ws.Cells.Item(36, 1).Value = "FY_4.png"
$ws.Cells.Item(36, 2).Value = 0.993
$ws.Cells.Item(36, 3).Value = 0.003
$ws.Cells.Item(36, 4).Value = 0.383
$ws.Cells.Item(36, 5).Value = 0.6
$ws.Cells.Item(36, 6).Value = "Fanny Yusuf"
$ws.Cells.Item(36, 7).Value = "Benar"
$ws.Cells.Item(37, 2).Value = 0.667
$ws.Cells.Item(37, 4).Value = 0.397
$ws.Cells.Item(37, 5).Value = 0.8
$ws.Cells.Item(38, 2).Value = 0.889
$ws.Cells.Item(38, 3).Value = 0.003
$ws.Cells.Item(38, 4).Value = 0.392
$ws.Cells.Item(39, 2).Value = 0.736
$ws.Cells.Item(39, 4).Value = 0.391
$ws.Cells.Item(40, 2).Value = 2.084
$ws.Cells.Item(40, 3).Value = 0.007
$ws.Cells.Item(40, 4).Value = 0.099
$ws.Cells.Item(40, 5).Value = 0.9330000000000001
$ws.Cells.Item(40, 6).Value = "Tiara Oktavian"
$ws.Cells.Item(40, 7).Value = "Benar"
$ws.Cells.Item(41, 2).Value = 1.809
$ws.Cells.Item(41, 3).Value = 0.006
$ws.Cells.Item(41, 4).Value = 0.08400000000000001
$ws.Cells.Item(41, 5).Value = 0.9330000000000001
$ws.Cells.Item(42, 2).Value = 1.555
$ws.Cells.Item(42, 3).Value = 0.005
$ws.Cells.Item(42, 4).Value = 0.028
$ws.Cells.Item(42, 5).Value = 0.333
$ws.Cells.Item(43, 2).Value = 1.668
$ws.Cells.Item(43, 3).Value = 0.005
$ws.Cells.Item(43, 4).Value = 0.045
$ws.Cells.Item(44, 2).Value = 0.894
$ws.Cells.Item(44, 3).Value = 0.003
$ws.Cells.Item(44, 4).Value = 0.196
$ws.Cells.Item(45, 2).Value = 0.905
$ws.Cells.Item(45, 3).Value = 0.003
$ws.Cells.Item(45, 4).Value = 0.11
$ws.Cells.Item(45, 5).Value = 0.267
$ws.Cells.Item(45, 6).Value = "Tidak Diketahui"
$ws.Cells.Item(45, 7).Value = "Benar"
